$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.675.77"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.742.00"
$ws.Range("E3").Value = "  +8.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.20"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.59"
$ws.Range("E6").Value = "  -4.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.740.76"
$ws.Range("E7").Value = "  +8.74%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.36"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.95"
$ws.Range("E13").Value = "  +6.86%  "
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.361.54"
$ws.Range("E15").Value = "  +8.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.739.59"
$ws.Range("E16").Value = "  +8.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.772.83"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("E19").Value = "  +3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.12"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.77"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.56"
$ws.Range("E22").Value = "  +9.38%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.43"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.07"
$ws.Range("E27").Value = "  +3.17%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +18.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.52"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("E32").Value = "  +5.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.26"
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.340"
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("E40").Value = "  +4.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.36"
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.83"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.37"
$ws.Range("E43").Value = "  -7.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "422.79"
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.096.15"
$ws.Range("E45").Value = "  +5.34%  "
$ws.Range("E46").Value = "  -4.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0365"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.93"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.54"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.30"
$ws.Range("E51").Value = "  +0.29%  "
